$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Header updates on sheet1 (LP1912) ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 16:52:27"
$ws1.Cells.Item(3,1).Value = "Total filas: 392"

# --- Data row updates on sheet1 (LP1912), rows 6-397 ---
$ws1.Cells.Item(52,1).Value = "05:51:32"
$ws1.Cells.Item(52,2).Value = "07:31"
$ws1.Cells.Item(52,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(52,4).Value = 100
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "05:51:32"
$ws1.Cells.Item(53,2).Value = "07:31"
$ws1.Cells.Item(53,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(53,4).Value = 100
$ws1.Cells.Item(53,5).Value = "LP1912"
$ws1.Cells.Item(76,1).Value = "07:46:15"
$ws1.Cells.Item(76,2).Value = "08:33"
$ws1.Cells.Item(76,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(76,4).Value = 47
$ws1.Cells.Item(76,5).Value = "LP1912"
$ws1.Cells.Item(77,1).Value = "07:46:15"
$ws1.Cells.Item(77,2).Value = "08:33"
$ws1.Cells.Item(77,3).Value = "10_OLMOS"
$ws1.Cells.Item(77,4).Value = 47
$ws1.Cells.Item(77,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "09:38:04"
$ws1.Cells.Item(124,2).Value = "09:41"
$ws1.Cells.Item(124,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(124,4).Value = 3
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "09:38:04"
$ws1.Cells.Item(125,2).Value = "09:41"
$ws1.Cells.Item(125,3).Value = "14_ABASTO"
$ws1.Cells.Item(125,4).Value = 3
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "08:39:56"
$ws1.Cells.Item(126,2).Value = "09:41"
$ws1.Cells.Item(126,3).Value = "215C_EL PATO"
$ws1.Cells.Item(126,4).Value = 62
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(174,1).Value = "11:23:54"
$ws1.Cells.Item(174,2).Value = "11:25"
$ws1.Cells.Item(174,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(174,4).Value = 2
$ws1.Cells.Item(174,5).Value = "LP1912"
$ws1.Cells.Item(175,1).Value = "09:38:04"
$ws1.Cells.Item(175,2).Value = "11:25"
$ws1.Cells.Item(175,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(175,4).Value = 107
$ws1.Cells.Item(175,5).Value = "LP1912"
$ws1.Cells.Item(199,1).Value = "10:57:58"
$ws1.Cells.Item(199,2).Value = "12:06"
$ws1.Cells.Item(199,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(199,4).Value = 69
$ws1.Cells.Item(199,5).Value = "LP1912"
$ws1.Cells.Item(200,1).Value = "10:57:58"
$ws1.Cells.Item(200,2).Value = "12:06"
$ws1.Cells.Item(200,3).Value = "14_ABASTO"
$ws1.Cells.Item(200,4).Value = 69
$ws1.Cells.Item(200,5).Value = "LP1912"
$ws1.Cells.Item(201,1).Value = "10:28:12"
$ws1.Cells.Item(201,2).Value = "12:06"
$ws1.Cells.Item(201,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(201,4).Value = 98
$ws1.Cells.Item(201,5).Value = "LP1912"
$ws1.Cells.Item(212,1).Value = "10:57:58"
$ws1.Cells.Item(212,2).Value = "12:20"
$ws1.Cells.Item(212,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(212,4).Value = 83
$ws1.Cells.Item(212,5).Value = "LP1912"
$ws1.Cells.Item(213,1).Value = "10:57:58"
$ws1.Cells.Item(213,2).Value = "12:20"
$ws1.Cells.Item(213,3).Value = "215A_EL PATO"
$ws1.Cells.Item(213,4).Value = 83
$ws1.Cells.Item(213,5).Value = "LP1912"
$ws1.Cells.Item(224,1).Value = "11:51:05"
$ws1.Cells.Item(224,2).Value = "12:37"
$ws1.Cells.Item(224,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(224,4).Value = 46
$ws1.Cells.Item(224,5).Value = "LP1912"
$ws1.Cells.Item(226,1).Value = "11:51:05"
$ws1.Cells.Item(226,2).Value = "12:37"
$ws1.Cells.Item(226,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(226,4).Value = 46
$ws1.Cells.Item(226,5).Value = "LP1912"
$ws1.Cells.Item(249,1).Value = "11:23:54"
$ws1.Cells.Item(249,2).Value = "13:20"
$ws1.Cells.Item(249,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(249,4).Value = 117
$ws1.Cells.Item(249,5).Value = "LP1912"
$ws1.Cells.Item(250,1).Value = "11:51:05"
$ws1.Cells.Item(250,2).Value = "13:20"
$ws1.Cells.Item(250,3).Value = "10_OLMOS"
$ws1.Cells.Item(250,4).Value = 89
$ws1.Cells.Item(250,5).Value = "LP1912"
$ws1.Cells.Item(298,1).Value = "14:16:51"
$ws1.Cells.Item(298,2).Value = "15:04"
$ws1.Cells.Item(298,3).Value = "10_OLMOS"
$ws1.Cells.Item(298,4).Value = 48
$ws1.Cells.Item(298,5).Value = "LP1912"
$ws1.Cells.Item(299,1).Value = "14:40:41"
$ws1.Cells.Item(299,2).Value = "15:04"
$ws1.Cells.Item(299,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(299,4).Value = 24
$ws1.Cells.Item(299,5).Value = "LP1912"
$ws1.Cells.Item(309,1).Value = "13:30:15"
$ws1.Cells.Item(309,2).Value = "15:20"
$ws1.Cells.Item(309,3).Value = "15_ABASTO"
$ws1.Cells.Item(309,4).Value = 110
$ws1.Cells.Item(309,5).Value = "LP1912"
$ws1.Cells.Item(310,1).Value = "13:30:15"
$ws1.Cells.Item(310,2).Value = "15:20"
$ws1.Cells.Item(310,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(310,4).Value = 110
$ws1.Cells.Item(310,5).Value = "LP1912"
$ws1.Cells.Item(360,1).Value = "16:14:52"
$ws1.Cells.Item(360,2).Value = "16:56"
$ws1.Cells.Item(360,3).Value = "10_OLMOS"
$ws1.Cells.Item(360,4).Value = 42
$ws1.Cells.Item(360,5).Value = "LP1912"
$ws1.Cells.Item(361,1).Value = "15:19:52"
$ws1.Cells.Item(361,2).Value = "16:56"
$ws1.Cells.Item(361,3).Value = "17_179 Y 38"
$ws1.Cells.Item(361,4).Value = 97
$ws1.Cells.Item(361,5).Value = "LP1912"
$ws1.Cells.Item(377,1).Value = "16:52:27"
$ws1.Cells.Item(377,2).Value = "17:34"
$ws1.Cells.Item(377,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(377,4).Value = 42
$ws1.Cells.Item(377,5).Value = "LP1912"
$ws1.Cells.Item(378,1).Value = "16:52:27"
$ws1.Cells.Item(378,2).Value = "17:34"
$ws1.Cells.Item(378,3).Value = "10_OLMOS"
$ws1.Cells.Item(378,4).Value = 42
$ws1.Cells.Item(378,5).Value = "LP1912"
$ws1.Cells.Item(379,1).Value = "16:32:38"
$ws1.Cells.Item(379,2).Value = "17:35"
$ws1.Cells.Item(379,3).Value = "15_ABASTO"
$ws1.Cells.Item(379,4).Value = 63
$ws1.Cells.Item(379,5).Value = "LP1912"
$ws1.Cells.Item(380,1).Value = "16:45:22"
$ws1.Cells.Item(380,2).Value = "17:35"
$ws1.Cells.Item(380,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(380,4).Value = 50
$ws1.Cells.Item(380,5).Value = "LP1912"
$ws1.Cells.Item(381,1).Value = "16:14:52"
$ws1.Cells.Item(381,2).Value = "17:35"
$ws1.Cells.Item(381,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(381,4).Value = 81
$ws1.Cells.Item(381,5).Value = "LP1912"
$ws1.Cells.Item(382,1).Value = "15:51:40"
$ws1.Cells.Item(382,2).Value = "17:36"
$ws1.Cells.Item(382,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(382,4).Value = 105
$ws1.Cells.Item(382,5).Value = "LP1912"
$ws1.Cells.Item(383,1).Value = "16:45:22"
$ws1.Cells.Item(383,2).Value = "17:37"
$ws1.Cells.Item(383,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(383,4).Value = 52
$ws1.Cells.Item(383,5).Value = "LP1912"
$ws1.Cells.Item(384,1).Value = "15:51:40"
$ws1.Cells.Item(384,2).Value = "17:38"
$ws1.Cells.Item(384,3).Value = "17_ROMERO"
$ws1.Cells.Item(384,4).Value = 107
$ws1.Cells.Item(384,5).Value = "LP1912"
$ws1.Cells.Item(385,1).Value = "15:51:40"
$ws1.Cells.Item(385,2).Value = "17:40"
$ws1.Cells.Item(385,3).Value = "215B_EL PATO"
$ws1.Cells.Item(385,4).Value = 109
$ws1.Cells.Item(385,5).Value = "LP1912"
$ws1.Cells.Item(386,1).Value = "16:52:27"
$ws1.Cells.Item(386,2).Value = "17:40"
$ws1.Cells.Item(386,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(386,4).Value = 48
$ws1.Cells.Item(386,5).Value = "LP1912"
$ws1.Cells.Item(387,1).Value = "16:45:22"
$ws1.Cells.Item(387,2).Value = "17:41"
$ws1.Cells.Item(387,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(387,4).Value = 56
$ws1.Cells.Item(387,5).Value = "LP1912"
$ws1.Cells.Item(388,1).Value = "16:45:22"
$ws1.Cells.Item(388,2).Value = "17:45"
$ws1.Cells.Item(388,3).Value = "15_ABASTO"
$ws1.Cells.Item(388,4).Value = 60
$ws1.Cells.Item(388,5).Value = "LP1912"
$ws1.Cells.Item(389,1).Value = "15:51:40"
$ws1.Cells.Item(389,2).Value = "17:50"
$ws1.Cells.Item(389,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(389,4).Value = 119
$ws1.Cells.Item(389,5).Value = "LP1912"
$ws1.Cells.Item(390,1).Value = "16:14:52"
$ws1.Cells.Item(390,2).Value = "17:52"
$ws1.Cells.Item(390,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(390,4).Value = 98
$ws1.Cells.Item(390,5).Value = "LP1912"
$ws1.Cells.Item(391,1).Value = "16:14:52"
$ws1.Cells.Item(391,2).Value = "18:04"
$ws1.Cells.Item(391,3).Value = "17_ROMERO"
$ws1.Cells.Item(391,4).Value = 110
$ws1.Cells.Item(391,5).Value = "LP1912"
$ws1.Cells.Item(392,1).Value = "16:52:27"
$ws1.Cells.Item(392,2).Value = "18:08"
$ws1.Cells.Item(392,3).Value = "14_ABASTO"
$ws1.Cells.Item(392,4).Value = 76
$ws1.Cells.Item(392,5).Value = "LP1912"
$ws1.Cells.Item(393,1).Value = "16:32:38"
$ws1.Cells.Item(393,2).Value = "18:21"
$ws1.Cells.Item(393,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(393,4).Value = 109
$ws1.Cells.Item(393,5).Value = "LP1912"
$ws1.Cells.Item(394,1).Value = "16:32:38"
$ws1.Cells.Item(394,2).Value = "18:27"
$ws1.Cells.Item(394,3).Value = "215C_EL PATO"
$ws1.Cells.Item(394,4).Value = 115
$ws1.Cells.Item(394,5).Value = "LP1912"
$ws1.Cells.Item(395,1).Value = "16:45:22"
$ws1.Cells.Item(395,2).Value = "18:28"
$ws1.Cells.Item(395,3).Value = "215C_EL PATO"
$ws1.Cells.Item(395,4).Value = 103
$ws1.Cells.Item(395,5).Value = "LP1912"
$ws1.Cells.Item(396,1).Value = "16:45:22"
$ws1.Cells.Item(396,2).Value = "18:32"
$ws1.Cells.Item(396,3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(396,4).Value = 107
$ws1.Cells.Item(396,5).Value = "LP1912"
$ws1.Cells.Item(397,1).Value = "16:52:27"
$ws1.Cells.Item(397,2).Value = "18:48"
$ws1.Cells.Item(397,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(397,4).Value = 116
$ws1.Cells.Item(397,5).Value = "LP1912"

# --- sheet2 (LP1912-215) and sheet3 (6203-6173): update "Última actualización" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,1).Value = "Última actualización: 16:52:27"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,1).Value = "Última actualización: 16:52:27"
